$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E2").Value = '[''Normal'']'
$ws.Range("D3").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E3").Value = '[''Normal'']'
$ws.Range("D11").Value = '[1, 0, 1, 0, 1, 0, 0]'
$ws.Range("E11").Value = '[''Normal'', ''HardwareFault'', ''RegulationViolation'']'
$ws.Range("D12").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E12").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D15").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E15").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D24").Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Range("E24").Value = '[''HardwareFault'']'
$ws.Range("D25").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E25").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D26").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E26").Value = '[''SoftwareFault'']'
$ws.Range("D27").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E27").Value = '[''SoftwareFault'']'
$ws.Range("D29").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E29").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D31").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E31").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D56").Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Range("E56").Value = '[]'
$ws.Range("D58").Value = '[0, 0, 0, 1, 0, 0, 1]'
$ws.Range("E58").Value = '[''ParamViolation'', ''SoftwareFault'']'
$ws.Range("D73").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E73").Value = '[''Normal'']'
$ws.Range("D74").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E74").Value = '[''Normal'']'
$ws.Range("D82").Value = '[1, 1, 1, 0, 0, 0, 0]'
$ws.Range("E82").Value = '[''Normal'', ''SurroundingEnvironment'', ''HardwareFault'']'
$ws.Range("D92").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E92").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D107").Value = '[1, 0, 0, 0, 0, 1, 1]'
$ws.Range("E107").Value = '[''Normal'', ''CommunicationIssue'', ''SoftwareFault'']'
$ws.Range("D109").Value = '[1, 1, 0, 0, 0, 0, 1]'
$ws.Range("E109").Value = '[''Normal'', ''SurroundingEnvironment'', ''SoftwareFault'']'
